# Update LR-pairs data (Il33-Il1rl1) with newly recomputed TPM stats.
# The "ECs" target-cluster rows (D = "ECs") are dropped entirely, leaving
# only the three "FAPs" target-cluster rows, whose metric columns are
# refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows whose Target cluster (column D) is "ECs" -- these were
# sheet rows 2, 4 and 6. Delete from the bottom up so earlier row numbers
# stay valid.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()

# Row 2 (was row 3): ECs -> FAPs
$ws.Range("G2").Value = 0.1004223333333333
$ws.Range("H2").Value = 0.301267
$ws.Range("I2").Value = 0.002425263903734378
$ws.Range("J2").Value = 0.002425263903734379
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.3262840108353333
$ws.Range("R2").Value = 2.936556097518
$ws.Range("S2").Value = 0.002425263903734378
$ws.Range("T2").Value = 0.002425263903734379

# Row 3 (was row 5): FAPs -> FAPs
$ws.Range("I3").Value = 0.9126643201687427
$ws.Range("J3").Value = 0.9126643201687429
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("S3").Value = 0.9126643201687427
$ws.Range("T3").Value = 0.9126643201687429

# Row 4 (was row 7): MuSCs -> FAPs
$ws.Range("G4").Value = 3.515865666666667
$ws.Range("H4").Value = 10.547597
$ws.Range("I4").Value = 0.08491041592752281
$ws.Range("J4").Value = 0.08491041592752283
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 11.42346242314867
$ws.Range("R4").Value = 102.811161808338
$ws.Range("S4").Value = 0.08491041592752281
$ws.Range("T4").Value = 0.08491041592752283
